$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "vat" column (column L) header and values
$ws.Range("L1").Value = "vat"
$ws.Range("L2").Value = 5000
$ws.Range("L3").Value = 2000
$ws.Range("L4").Value = 1000
$ws.Range("L5").Value = 3000

# Copy the style from an existing header/data cell in the same row so the
# new column matches the rest of the table's formatting
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("K2:K5").Copy() | Out-Null
$ws.Range("L2:L5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update the active selection to match the final state
$ws.Range("L3").Select() | Out-Null
